$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 236 - this shifts the existing row 236
# (and everything below it, down to 311) down by one row to 237..312,
# growing the used range from A1:R311 to A1:R312.
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row 236 with the new weekly price record.
$ws.Cells.Item(236,1).Value = 4
$ws.Cells.Item(236,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(236,3).Value = "Los Lagos"
$ws.Cells.Item(236,4).Value = 44876
$ws.Cells.Item(236,5).Value = 10
$ws.Cells.Item(236,6).Value = 100112032
$ws.Cells.Item(236,7).Value = "Zapallo italiano"
$ws.Cells.Item(236,8).Value = "Sin especificar"
$ws.Cells.Item(236,9).Value = "Primera"
$ws.Cells.Item(236,10).Value = 250
$ws.Cells.Item(236,11).Value = 15000
$ws.Cells.Item(236,12).Value = 15000
$ws.Cells.Item(236,13).Value = 15000
$ws.Cells.Item(236,14).Value = "$/caja 50 unidades"
$ws.Cells.Item(236,15).Value = "Región de O'Higgins"
$ws.Cells.Item(236,16).Value = 300
$ws.Cells.Item(236,17).Value = 50
$ws.Cells.Item(236,18).Value = "Hortaliza"
